$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 5..18 mapped to column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 16
$ws1.Range("F6").Value = 14
$ws1.Range("F7").Value = 557
$ws1.Range("F8").Value = 7813
$ws1.Range("F9").Value = 747
$ws1.Range("F10").Value = 216
$ws1.Range("F11").Value = 1092
$ws1.Range("F12").Value = 716
$ws1.Range("F15").Value = 191
$ws1.Range("F16").Value = 38
$ws1.Range("F17").Value = 204
$ws1.Range("F18").Value = 798

# Sheet "全部类型" (All types) - rows 5..19 mapped to column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 16
$ws4.Range("F6").Value = 14
$ws4.Range("F8").Value = 557
$ws4.Range("F9").Value = 7813
$ws4.Range("F10").Value = 747
$ws4.Range("F11").Value = 216
$ws4.Range("F12").Value = 1092
$ws4.Range("F13").Value = 716
$ws4.Range("F16").Value = 191
$ws4.Range("F17").Value = 38
$ws4.Range("F18").Value = 204
$ws4.Range("F19").Value = 798
